$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format on the Price column data range so that
# numeric-looking strings (e.g. "0.693") are stored as text, matching
# the original inline-string cell type, then clear the format again
# so no residual style index is left on the cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '43.869.17'
$ws.Range("E2").Value = '  +0.89%  '

$ws.Range("D3").Value = '2.354.68'
$ws.Range("E3").Value = '  -0.54%  '

$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").Value = '0.693'
$ws.Range("E5").Value = '  +6.43%  '

$ws.Range("D6").Value = '240.94'
$ws.Range("E6").Value = '  +2.93%  '

$ws.Range("D7").Value = '76.79'
$ws.Range("E7").Value = '  +5.69%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("D9").Value = '0.632'
$ws.Range("E9").Value = '  +24.79%  '

$ws.Range("D10").Value = '0.102'
$ws.Range("E10").Value = '  +3.94%  '

$ws.Range("D11").Value = '57.39'
$ws.Range("E11").Value = '  +0.69%  '

$ws.Range("D12").Value = '34.14'
$ws.Range("E12").Value = '  +24.38%  '

$ws.Range("E13").Value = '  +19.53%  '

$ws.Range("D15").Value = '2.703.57'
$ws.Range("E15").Value = '  -0.39%  '

$ws.Range("D16").Value = '16.89'
$ws.Range("E16").Value = '  +2.95%  '

$ws.Range("E17").Value = '  +5.86%  '

$ws.Range("D18").Value = '2.352.34'
$ws.Range("E18").Value = '  -0.49%  '

$ws.Range("D19").Value = '43.820.12'
$ws.Range("E19").Value = '  +0.81%  '

$ws.Range("D20").Value = '0.0000103'
$ws.Range("E20").Value = '  +2.58%  '

$ws.Range("D21").Value = '6.65'
$ws.Range("E21").Value = '  +4.33%  '

$ws.Range("D22").Value = '77.59'
$ws.Range("E22").Value = '  +2.79%  '

$ws.Range("D23").Value = '255.88'
$ws.Range("E23").Value = '  +1.70%  '

$ws.Range("E25").Value = '  +2.81%  '

$ws.Range("D26").Value = '11.04'
$ws.Range("E26").Value = '  +8.80%  '

$ws.Range("E27").Value = '  -5.10%  '

$ws.Range("E28").Value = '  +15.41%  '

$ws.Range("E29").Value = '  +2.20%  '

$ws.Range("D30").Value = '23.04'
$ws.Range("E30").Value = '  +1.56%  '

$ws.Range("D31").Value = '174.90'
$ws.Range("E31").Value = '  +1.21%  '

$ws.Range("E32").Value = '  -3.25%  '

$ws.Range("E33").Value = '  +6.18%  '

$ws.Range("D34").Value = '0.0759'
$ws.Range("E34").Value = '  +8.48%  '

$ws.Range("D35").Value = '5.31'
$ws.Range("E35").Value = '  +5.35%  '

$ws.Range("D36").Value = '5.40'
$ws.Range("E36").Value = '  +6.43%  '

$ws.Range("D37").Value = '3.81'
$ws.Range("E37").Value = '  +1.23%  '

$ws.Range("D38").Value = '2.43'
$ws.Range("E38").Value = '  -0.59%  '

$ws.Range("D39").Value = '6.47'
$ws.Range("E39").Value = '  -2.78%  '

$ws.Range("E40").Value = '  +8.17%  '

$ws.Range("D41").Value = '19.40'
$ws.Range("E41").Value = '  -0.62%  '

$ws.Range("E42").Value = '  +17.39%  '

$ws.Range("D43").Value = '9.01'
$ws.Range("E43").Value = '  +0.79%  '

$ws.Range("E44").Value = '  +0.11%  '

$ws.Range("E45").Value = '  +6.67%  '

$ws.Range("E46").Value = '  +14.16%  '

$ws.Range("E47").Value = '  +4.29%  '

$ws.Range("D48").Value = '101.88'
$ws.Range("E48").Value = '  +1.69%  '

$ws.Range("E49").Value = '  +1.79%  '

$ws.Range("D50").Value = '4.50'
$ws.Range("E50").Value = '  -0.72%  '

$ws.Range("D51").Value = '54.96'
$ws.Range("E51").Value = '  +7.77%  '

$ws.Range("D2:D51").ClearFormats()
